$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 61.333332
$ws.Range("I8").Value = 61.333332
$ws.Range("K8").Value = 183.999996
$ws.Range("M8").Value = -44.99999600000001

$ws.Range("H11").Value = 89038.5
$ws.Range("I11").Value = 89038.5
$ws.Range("K11").Value = 89038.5
$ws.Range("M11").Value = -88898.5

$ws.Range("H19").Value = 1433.375
$ws.Range("J19").Value = 1666.6666
$ws.Range("L19").Value = 1666.6666
$ws.Range("N19").Value = -2016.6666

$ws.Range("H32").Value = 1000
$ws.Range("J32").Value = 1000
$ws.Range("L32").Value = 1000
$ws.Range("N32").Value = -1652

$ws.Range("H39").Value = 483.83334
$ws.Range("I39").Value = 301
$ws.Range("K39").Value = 903
$ws.Range("M39").Value = -607

$ws.Range("H53").Value = 463.25
$ws.Range("I53").Value = 499
$ws.Range("J53").Value = 458.14285
$ws.Range("K53").Value = 499
$ws.Range("L53").Value = 458.14285
$ws.Range("M53").Value = 138
$ws.Range("N53").Value = -1732.14285

$ws.Range("H55").Value = 230
$ws.Range("J55").Value = 150
$ws.Range("L55").Value = 150
$ws.Range("N55").Value = -578

$ws.Range("H86").Value = 4000
$ws.Range("J86").Value = 4000
$ws.Range("L86").Value = 4000
$ws.Range("N86").Value = -6246

$ws.Range("H89").Value = 4000
$ws.Range("J89").Value = 4000
$ws.Range("L89").Value = 20000
$ws.Range("N89").Value = -31232

$ws.Range("H111").Value = 799.5
$ws.Range("J111").Value = 999
$ws.Range("L111").Value = 2997
$ws.Range("N111").Value = -9131

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 537.4286
$ws.Range("I80").Value = 57
$ws.Range("J80").Value = 617.5
$ws.Range("K80").Value = 57
$ws.Range("L80").Value = 617.5
$ws.Range("M80").Value = 941
$ws.Range("N80").Value = -2613.5

$ws.Range("H83").Value = 537.4286
$ws.Range("I83").Value = 57
$ws.Range("J83").Value = 617.5
$ws.Range("K83").Value = 285
$ws.Range("L83").Value = 3087.5
$ws.Range("M83").Value = 4707
$ws.Range("N83").Value = -13071.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3051.75
$ws.Range("I16").Value = 2404
$ws.Range("J16").Value = 4995
$ws.Range("K16").Value = 2404
$ws.Range("L16").Value = 4995
$ws.Range("M16").Value = -2117
$ws.Range("N16").Value = -5569

$ws.Range("H62").Value = 20501.5
$ws.Range("J62").Value = 20501.5
$ws.Range("L62").Value = 20501.5
$ws.Range("N62").Value = -21749.5

$ws.Range("H65").Value = 20501.5
$ws.Range("J65").Value = 20501.5
$ws.Range("L65").Value = 102507.5
$ws.Range("N65").Value = -108747.5

$ws.Range("H86").Value = 20426.715
$ws.Range("I86").Value = 20426.715
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 20426.715
$ws.Range("L86").Value = 0
$ws.Range("N86").Value = -19303.715
$ws.Range("M86").ClearContents()

$ws.Range("H89").Value = 20426.715
$ws.Range("I89").Value = 20426.715
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 102133.575
$ws.Range("L89").Value = 0
$ws.Range("N89").Value = -96517.575
$ws.Range("M89").ClearContents()

$ws.Range("H107").Value = 1778.5
$ws.Range("I107").Value = 1443.25
$ws.Range("J107").Value = 2449
$ws.Range("K107").Value = 1443.25
$ws.Range("L107").Value = 2449
$ws.Range("M107").Value = 476.75
$ws.Range("N107").Value = -6289

$ws.Range("H113").Value = 3051.75
$ws.Range("I113").Value = 2404
$ws.Range("J113").Value = 4995
$ws.Range("K113").Value = 2404
$ws.Range("L113").Value = 4995
$ws.Range("M113").Value = -234
$ws.Range("N113").Value = -9335

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 262
$ws.Range("J23").Value = 500
$ws.Range("L23").Value = 1500
$ws.Range("N23").Value = -1970

$ws.Range("H29").Value = 1142.5
$ws.Range("I29").Value = 280
$ws.Range("J29").Value = 1315
$ws.Range("K29").Value = 840
$ws.Range("L29").Value = 3945
$ws.Range("M29").Value = -563
$ws.Range("N29").Value = -4499

$ws.Range("H40").Value = 433.66666
$ws.Range("I40").Value = 501
$ws.Range("J40").Value = 400
$ws.Range("K40").Value = 2004
$ws.Range("L40").Value = 1600
$ws.Range("M40").Value = -1935
$ws.Range("N40").Value = -1738

$ws.Range("H92").Value = 438
$ws.Range("I92").Value = 362.33334
$ws.Range("K92").Value = 1087.00002
$ws.Range("M92").Value = 160.9999800000001

$ws.Range("H105").Value = 10000
$ws.Range("J105").Value = 10000
$ws.Range("L105").Value = 30000
$ws.Range("N105").Value = -35242

$ws.Range("H134").Value = 3096.25
$ws.Range("I134").Value = 3076.6
$ws.Range("K134").Value = 9229.799999999999
$ws.Range("M134").Value = -4159.799999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").ClearContents()

$ws.Range("H16").Value = 1312.25
$ws.Range("I16").Value = 1416.3334
$ws.Range("K16").Value = 1416.3334
$ws.Range("M16").Value = -1246.3334

$ws.Range("H22").Value = 9829.058999999999
$ws.Range("I22").Value = 12049.375
$ws.Range("J22").Value = 7855.4443
$ws.Range("K22").Value = 12049.375
$ws.Range("L22").Value = 7855.4443
$ws.Range("M22").Value = -11754.375
$ws.Range("N22").Value = -8445.444299999999

$ws.Range("H27").Value = 9829.058999999999
$ws.Range("I27").Value = 12049.375
$ws.Range("J27").Value = 7855.4443
$ws.Range("K27").Value = 12049.375
$ws.Range("L27").Value = 7855.4443
$ws.Range("M27").Value = -11942.375
$ws.Range("N27").Value = -8069.4443

$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()

$ws.Range("H132").Value = 5000
$ws.Range("I132").Value = 5000
$ws.Range("K132").Value = 15000
$ws.Range("M132").Value = -12470

$ws.Range("H136").Value = 1837.5
$ws.Range("I136").Value = 1837.5
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 5512.5
$ws.Range("L136").Value = 0
$ws.Range("N136").Value = -2962.5
$ws.Range("M136").ClearContents()
